# Apply "Fixing network data cleanining scripts" edits to NEW_HAMPSHIRE_2021 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns to clean/English machine-readable names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Normalize capitalization of "de"/"los" -> "De"/"Los" in place names
$ws.Range("A14").Value = "Ciudad De México"
$ws.Range("A20").Value = "Estado De México"
$ws.Range("B26").Value = "Acapulco De Juárez"
$ws.Range("B31").Value = "Mártir De Cuilapan"
$ws.Range("B36").Value = "Progreso De Obregón"
$ws.Range("B38").Value = "Autlán De Navarro"
$ws.Range("B39").Value = "Cuautitlán De García Barragán"
$ws.Range("B41").Value = "Encarnación De Díaz"
$ws.Range("B45").Value = "Unión De Tula"
$ws.Range("B55").Value = "San Nicolás De Los Garza"
$ws.Range("B57").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B69").Value = "Jalpan De Serra"
$ws.Range("B70").Value = "Landa De Matamoros"

# 3. Tiny floating point correction (re-computed percentage)
$ws.Range("D72").Value = 0.09352517985611512

# 4. Remove trailing footnote/source rows (97-101); data now ends at row 95
$ws.Rows("97:101").Delete()
